$d = $word.ActiveDocument

# Locate the (single) field that holds the m:'...'.representationByName()....
# instrText token stream we need to rewrite as plain literal text runs
# wrapped in "{ ... }" (mirrors TokenIteratorFieldRewriterSplit behaviour).
$field = $d.Fields(1)

$fieldStart = $field.Code.Start - 2   # back up over the begin fldChar run
$fieldEnd   = $field.Result.End + 1   # past the end fldChar run

$rng = $d.Range($fieldStart, $fieldEnd)

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$texts = @(
    "{m:",
    "'",
    "anydsl class diagram",
    "'.",
    "r",
    "epresentation",
    "By",
    "Name",
    "()",
    ".",
    "asImage(",
    "'JP",
    "E",
    "G'",
    ")",
    ".fit(400, 400)}"
)

$runsXml = ""
foreach ($t in $texts) {
    $runsXml += "<w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>$t</w:t></w:r>"
}

$xml = "<w:p xmlns:w=`"$w`" xmlns:w14=`"http://schemas.microsoft.com/office/word/2010/wordml`" w14:paraId=`"70197535`" w14:textId=`"17BBC239`" w:rsidR=`"00A10D75`" w:rsidRDefault=`"00474E78`" w:rsidP=`"00F65375`"><w:pPr><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr>$runsXml</w:p>"

$rng.InsertXML($xml)
